$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 130667093
$ws.Range("B2").Value = 57884
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = 'Tretåig hackspett'
$ws.Range("G2").Value = 'Picoides tridactylus'
$ws.Range("H2").Value = '(Linnaeus, 1758)'
$ws.Range("M2").Value = 'färska spår'
$ws.Range("Q2").Value = 491408
$ws.Range("R2").Value = 6759381

# Row 3
$ws.Range("A3").Value = 130670771
$ws.Range("B3").Value = 79243
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = 'Garnlav'
$ws.Range("G3").Value = 'Alectoria sarmentosa'
$ws.Range("H3").Value = '(Ach.) Ach.'
$ws.Range("M3").Value = ""
$ws.Range("Q3").Value = 491374
$ws.Range("R3").Value = 6759416

# Row 4
$ws.Range("A4").Value = 130670627
$ws.Range("B4").Value = 79243
$ws.Range("D4").Value = 'NT'
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = 'Garnlav'
$ws.Range("G4").Value = 'Alectoria sarmentosa'
$ws.Range("H4").Value = '(Ach.) Ach.'
$ws.Range("M4").Value = ""
$ws.Range("Q4").Value = 491376
$ws.Range("R4").Value = 6759442

# Row 11
$ws.Range("A11").Value = 130662234
$ws.Range("B11").Value = 8451
$ws.Range("D11").Value = 'LC'
$ws.Range("E11").Value = 106545
$ws.Range("F11").Value = 'Mindre märgborre'
$ws.Range("G11").Value = 'Tomicus minor'
$ws.Range("H11").Value = '(Hartig, 1834)'
$ws.Range("M11").Value = 'färska gnagspår'
$ws.Range("Q11").Value = 491455
$ws.Range("R11").Value = 6759425

# Row 12
$ws.Range("A12").Value = 130668644
$ws.Range("B12").Value = 79243
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = 'Garnlav'
$ws.Range("G12").Value = 'Alectoria sarmentosa'
$ws.Range("H12").Value = '(Ach.) Ach.'
$ws.Range("M12").Value = ""
$ws.Range("Q12").Value = 491424
$ws.Range("R12").Value = 6759256

# Row 14
$ws.Range("A14").Value = 130665997
$ws.Range("B14").Value = 8451
$ws.Range("D14").Value = 'LC'
$ws.Range("E14").Value = 106545
$ws.Range("F14").Value = 'Mindre märgborre'
$ws.Range("G14").Value = 'Tomicus minor'
$ws.Range("H14").Value = '(Hartig, 1834)'
$ws.Range("M14").Value = 'färska gnagspår'
$ws.Range("Q14").Value = 491418
$ws.Range("R14").Value = 6759396

# Row 15
$ws.Range("A15").Value = 130668452
$ws.Range("B15").Value = 8451
$ws.Range("D15").Value = 'LC'
$ws.Range("E15").Value = 106545
$ws.Range("F15").Value = 'Mindre märgborre'
$ws.Range("G15").Value = 'Tomicus minor'
$ws.Range("H15").Value = '(Hartig, 1834)'
$ws.Range("M15").Value = 'äldre gnagspår'
$ws.Range("Q15").Value = 491385
$ws.Range("R15").Value = 6759286

# Row 16
$ws.Range("A16").Value = 130661548
$ws.Range("B16").Value = 8451
$ws.Range("D16").Value = 'LC'
$ws.Range("E16").Value = 106545
$ws.Range("F16").Value = 'Mindre märgborre'
$ws.Range("G16").Value = 'Tomicus minor'
$ws.Range("H16").Value = '(Hartig, 1834)'
$ws.Range("M16").Value = 'äldre gnagspår'
$ws.Range("Q16").Value = 491487
$ws.Range("R16").Value = 6759357

# Row 17
$ws.Range("A17").Value = 130661510
$ws.Range("B17").Value = 79243
$ws.Range("D17").Value = 'NT'
$ws.Range("E17").Value = 6425
$ws.Range("F17").Value = 'Garnlav'
$ws.Range("G17").Value = 'Alectoria sarmentosa'
$ws.Range("H17").Value = '(Ach.) Ach.'
$ws.Range("M17").Value = ""
$ws.Range("Q17").Value = 491504
$ws.Range("R17").Value = 6759336

# Row 20
$ws.Range("A20").Value = 130661581
$ws.Range("B20").Value = 79243
$ws.Range("D20").Value = 'NT'
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = 'Garnlav'
$ws.Range("G20").Value = 'Alectoria sarmentosa'
$ws.Range("H20").Value = '(Ach.) Ach.'
$ws.Range("M20").Value = ""
$ws.Range("Q20").Value = 491481
$ws.Range("R20").Value = 6759380

# Row 21
$ws.Range("A21").Value = 130670281
$ws.Range("B21").Value = 57884
$ws.Range("D21").Value = 'NT'
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = 'Tretåig hackspett'
$ws.Range("G21").Value = 'Picoides tridactylus'
$ws.Range("H21").Value = '(Linnaeus, 1758)'
$ws.Range("M21").Value = 'färska spår'
$ws.Range("Q21").Value = 491315
$ws.Range("R21").Value = 6759520

# Row 22
$ws.Range("A22").Value = 130661613
$ws.Range("B22").Value = 79243
$ws.Range("D22").Value = 'NT'
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = 'Garnlav'
$ws.Range("G22").Value = 'Alectoria sarmentosa'
$ws.Range("H22").Value = '(Ach.) Ach.'
$ws.Range("M22").Value = ""
$ws.Range("Q22").Value = 491477
$ws.Range("R22").Value = 6759416
